# Update automatico via Actualizar 06-15-2020 15-14-03
#
# Sheet "trabajo": refresh the "Fecha publicación" (column I) values for a
# few rows, and move the view/selection the way the author left it after
# editing (scrolled down to row 14, cell I16 selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2  (ID_Dato 1)  : Fecha publicación 06-04-2020 (was 04-06-2020)
$ws.Range("I2").Value = 43986
# Row 3  (ID_Dato 2)  : Fecha publicación 02-06-2020 (was 06-02-2020)
$ws.Range("I3").Value = 43984
# Row 15 (ID_Dato 14) : Fecha publicación 02-04-2020 (was 04-02-2020)
$ws.Range("I15").Value = 43923

# Scroll the window so row 14 / column E is the top-left visible cell,
# then select I16 (matches the saved view state in the workbook).
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("I16").Select()
